# Update "想去人数" (interested-count) / price figures on each sheet to match
# the newer site snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 503
$ws.Range("F6").Value = 939
$ws.Range("F9").Value = 2226
$ws.Range("F10").Value = 638
$ws.Range("F11").Value = 301
$ws.Range("F12").Value = 126
$ws.Range("F13").Value = 1117
$ws.Range("F15").Value = 2249
$ws.Range("F17").Value = 14037
$ws.Range("F19").Value = 1307
$ws.Range("F20").Value = 66
$ws.Range("F21").Value = 569
$ws.Range("F22").Value = 142
$ws.Range("F23").Value = 39
$ws.Range("F26").Value = 44
$ws.Range("F27").Value = 284
$ws.Range("F29").Value = 9

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = 266
$ws.Range("F18").Value = 30
$ws.Range("F22").Value = 3

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5733

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 503
$ws.Range("F8").Value = 939
$ws.Range("G11").Value = 266
$ws.Range("F12").Value = 2226
$ws.Range("F13").Value = 638
$ws.Range("F14").Value = 301
$ws.Range("F16").Value = 126
$ws.Range("F18").Value = 1117
$ws.Range("F23").Value = 2249
$ws.Range("F28").Value = 1307
$ws.Range("F29").Value = 66
$ws.Range("F30").Value = 569
$ws.Range("F31").Value = 142
$ws.Range("F32").Value = 39
$ws.Range("F36").Value = 44
$ws.Range("F39").Value = 284
$ws.Range("F43").Value = 3
